$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rangeAddr, $value) {
    # Force the cell to remain a text cell (matches the workbook's original
    # inline-string cells) instead of letting Excel auto-coerce
    # numeric-looking strings (e.g. "306.77") into a Number cell.
    $rng = $ws.Range($rangeAddr)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue "D2" "42.940.85"
$ws.Range("E2").Value = "  -1.35%  "

# Row 3 - Ethereum
Set-TextValue "D3" "2.340.67"
$ws.Range("E3").Value = "  +1.24%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.02%  "

# Row 5 - BNB
Set-TextValue "D5" "306.77"
$ws.Range("E5").Value = "  -1.41%  "

# Row 6 - Solana
Set-TextValue "D6" "100.87"
$ws.Range("E6").Value = "  -1.37%  "

# Row 7 - XRP
$ws.Range("E7").Value = "  -4.87%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  +0.03%  "

# Row 9 - Cardano
$ws.Range("E9").Value = "  -3.68%  "

# Row 10 - Avalanche
Set-TextValue "D10" "34.95"
$ws.Range("E10").Value = "  -2.14%  "

# Row 11 - OKB
Set-TextValue "D11" "52.05"
$ws.Range("E11").Value = "  +0.19%  "

# Row 12 - Dogecoin
$ws.Range("E12").Value = "  -1.98%  "

# Row 13 - TRON
$ws.Range("E13").Value = "  -0.37%  "

# Row 14 - Polkadot
Set-TextValue "D14" "6.81"
$ws.Range("E14").Value = "  -3.04%  "

# Row 15 - Chainlink
Set-TextValue "D15" "15.87"
$ws.Range("E15").Value = "  +5.80%  "

# Row 16 - now WrappedEther (was Polygon)
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue "D16" "2.373.70"
$ws.Range("E16").Value = "  +2.43%  "

# Row 17 - now Polygon (was WrappedEther)
$ws.Range("B17").Value = "Polygon"
$ws.Range("C17").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextValue "D17" "0.808"
$ws.Range("E17").Value = "  -0.12%  "

# Row 18 - WrappedBTC
Set-TextValue "D18" "42.844.72"
$ws.Range("E18").Value = "  -1.33%  "

# Row 19 - Uniswap
Set-TextValue "D19" "6.23"
$ws.Range("E19").Value = "  +0.75%  "

# Row 20 - ShibaInu
$ws.Range("E20").Value = "  -1.59%  "

# Row 21 - InternetComputer(DFINITY)
Set-TextValue "D21" "11.68"
$ws.Range("E21").Value = "  -5.95%  "

# Row 22 - Litecoin
Set-TextValue "D22" "67.83"

# Row 23 - BitcoinCash
Set-TextValue "D23" "236.95"
$ws.Range("E23").Value = "  -1.90%  "

# Row 24 - ImmutableX
Set-TextValue "D24" "2.02"
$ws.Range("E24").Value = "  -0.88%  "

# Row 25 - PancakeSwap
Set-TextValue "D25" "2.56"
$ws.Range("E25").Value = "  -2.35%  "

# Row 26 - Dai
$ws.Range("E26").Value = "  -0.29%  "

# Row 27 - EthereumClassic
Set-TextValue "D27" "25.40"
$ws.Range("E27").Value = "  +2.64%  "

# Row 28 - Toncoin
$ws.Range("E28").Value = "  +3.14%  "

# Row 29 - InjectiveProtocol
Set-TextValue "D29" "35.09"
$ws.Range("E29").Value = "  -4.41%  "

# Row 30 - Cosmos
$ws.Range("E30").Value = "  -2.65%  "

# Row 31 - Monero
Set-TextValue "D31" "160.01"
$ws.Range("E31").Value = "  -4.70%  "

# Row 32 - FirstDigitalUSD
$ws.Range("E32").Value = "  -0.01%  "

# Row 33 - Filecoin
$ws.Range("E33").Value = "  -2.89%  "

# Row 34 - RenderToken
Set-TextValue "D34" "4.70"
$ws.Range("E34").Value = "  +9.68%  "

# Row 35 - WEMIXToken
Set-TextValue "D35" "2.48"
$ws.Range("E35").Value = "  -1.68%  "

# Row 37 - Hedera
$ws.Range("E37").Value = "  -2.27%  "

# Row 38 - LidoDAOToken
$ws.Range("E38").Value = "  -4.01%  "

# Row 40 - Kaspa
$ws.Range("E40").Value = "  -2.78%  "

# Row 41 - Stellar
$ws.Range("E41").Value = "  -2.55%  "

# Row 42 - Maker
Set-TextValue "D42" "2.018.30"
$ws.Range("E42").Value = "  +2.45%  "

# Row 43 - VeChain
$ws.Range("E43").Value = "  -1.12%  "

# Row 44 - EnergySwap
Set-TextValue "D44" "18.75"
$ws.Range("E44").Value = "  -2.91%  "

# Row 45 - FraxShare
Set-TextValue "D45" "10.29"
$ws.Range("E45").Value = "  +3.50%  "

# Row 46 - NEARProtocol
Set-TextValue "D46" "2.96"
$ws.Range("E46").Value = "  -0.76%  "

# Row 47 - MultiversX
Set-TextValue "D47" "56.38"
$ws.Range("E47").Value = "  +1.56%  "

# Row 48 - HuobiToken
$ws.Range("E48").Value = "  -0.73%  "

# Row 49 - RocketPoolETH
Set-TextValue "D49" "2.565.88"
$ws.Range("E49").Value = "  +1.05%  "

# Row 50 - THORChain
$ws.Range("E50").Value = "  +1.98%  "

# Row 51 - Stacks
Set-TextValue "D51" "1.51"
$ws.Range("E51").Value = "  -3.60%  "
